# Apply the vaccine price-list fixes described in the commit:
#   "jsut needs a few fixes. still parsing the last of the excel files for
#    issues. Almost perfect :)"
#
# Sheet 1 = "Pediatric VFC Vaccine "
# Sheet 2 = "Adult VFC Vaccine "
# Sheet 3 = "Pediatric influenza Influenza"

$wb = $excel.ActiveWorkbook

$wsPed   = $wb.Worksheets.Item(1)   # Pediatric VFC Vaccine
$wsAdult = $wb.Worksheets.Item(2)   # Adult VFC Vaccine
$wsFlu   = $wb.Worksheets.Item(3)   # Pediatric influenza Influenza

# --- Sheet 1: Pediatric VFC Vaccine ---------------------------------------

# Kinrix (rows 6-7) previously shared one combined packaging string for both
# rows; split it into the two real, distinct package descriptions.
$wsPed.Range("D6").Value = "10 pack - 1 dose vials "
$wsPed.Range("D7").Value = "5 pack - 1 dose T-L syringes "

# TriHIBit (row 11) was mislabeled with the generic "DTaP " vaccine name;
# correct it to "DTaP-Hib ".
$wsPed.Range("A11").Value = "DTaP-Hib "

# COMVAX (row 14) vaccine name corrected to include Hib.
$wsPed.Range("A14").Value = "Hepatitis B-Hib "

# ENGERIX B (rows 20-21) vaccine name: fix missing space/slash.
$wsPed.Range("A20").Value = "Hepatitis B Pediatric/Adolescent"
$wsPed.Range("A21").Value = "Hepatitis B Pediatric/Adolescent"

# RECOMBIVAX HB (row 22) vaccine name: fix missing slash (keeps double space).
$wsPed.Range("A22").Value = "Hepatitis B  Pediatric/Adolescent"

# Prevnar (row 29) vaccine name: add missing space.
$wsPed.Range("A29").Value = "Pneumococcal 7-valent (Pediatric)"

# --- Sheet 2: Adult VFC Vaccine -------------------------------------------

# DECAVAC (row 15) vaccine name had a typo duplicate string; fix spacing to
# match the canonical "Tetanus  Diphtheria Toxoids " text used elsewhere.
$wsAdult.Range("A15").Value = "Tetanus  Diphtheria Toxoids "

# --- Sheet 3: Pediatric influenza Influenza -------------------------------

$wsFlu.Range("A2").Value = "Influenza   (Age 6 months and older)"
$wsFlu.Range("A3").Value = "Influenza  (Age 6-35 months)"
$wsFlu.Range("A4").Value = "Influenza  (Age 36 months and older)"
$wsFlu.Range("A5").Value = "Influenza  (Age 36 months and older)"
$wsFlu.Range("A8").Value = "Influenza  Live, Intranasal (Age 2-49 years)"
